$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 85.5
$ws.Range("I9").Value = 85.71429000000001
$ws.Range("J9").Value = 85
$ws.Range("K9").Value = 85.71429000000001
$ws.Range("L9").Value = 85
$ws.Range("M9").Value = 83.28570999999999
$ws.Range("N9").Value = -423
$ws.Range("H17").Value = 1283.4634
$ws.Range("J17").Value = 1303.825
$ws.Range("L17").Value = 3911.475
$ws.Range("N17").Value = -4247.475
$ws.Range("H28").Value = 7262.5713
$ws.Range("I28").Value = 7262.5713
$ws.Range("K28").Value = 7262.5713
$ws.Range("M28").Value = -6777.5713
$ws.Range("H38").Value = 6977.074
$ws.Range("I38").Value = 6963.8823
$ws.Range("J38").Value = 6999.5
$ws.Range("K38").Value = 20891.6469
$ws.Range("L38").Value = 20998.5
$ws.Range("M38").Value = -20519.6469
$ws.Range("N38").Value = -21742.5
$ws.Range("H40").Value = 1199.6
$ws.Range("J40").Value = 899
$ws.Range("L40").Value = 899
$ws.Range("N40").Value = -1249
$ws.Range("H43").Value = 8127.3
$ws.Range("I43").Value = 4000
$ws.Range("J43").Value = 8585.888999999999
$ws.Range("K43").Value = 4000
$ws.Range("L43").Value = 8585.888999999999
$ws.Range("M43").Value = -3931
$ws.Range("N43").Value = -8723.888999999999
$ws.Range("H57").Value = 29272.727
$ws.Range("J57").Value = 29272.727
$ws.Range("L57").Value = 87818.181
$ws.Range("N57").Value = -88816.181
$ws.Range("H58").Value = 25005902
$ws.Range("I58").Value = 31250504
$ws.Range("K58").Value = 93751512
$ws.Range("M58").Value = -93751362
$ws.Range("H112").Value = 1711.9333
$ws.Range("J112").Value = 2124.125
$ws.Range("L112").Value = 6372.375
$ws.Range("N112").Value = -8588.375
$ws.Range("H135").Value = 1249.6154
$ws.Range("I135").Value = 1249.6154
$ws.Range("K135").Value = 11246.5386
$ws.Range("M135").Value = -8711.5386
$ws.Range("H138").Value = 2786.2554
$ws.Range("J138").Value = 3419.4546
$ws.Range("L138").Value = 10258.3638
$ws.Range("N138").Value = -20538.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 43750
$ws.Range("I109").Value = 30000
$ws.Range("J109").Value = 48333.332
$ws.Range("K109").Value = 30000
$ws.Range("L109").Value = 48333.332
$ws.Range("M109").Value = -28613
$ws.Range("N109").Value = -51107.332
$ws.Range("H110").Value = 1029.92
$ws.Range("I110").Value = 1029.92
$ws.Range("K110").Value = 1029.92
$ws.Range("M110").Value = 1015.08
$ws.Range("H122").Value = 2929.818
$ws.Range("I122").Value = 2998
$ws.Range("J122").Value = 2923
$ws.Range("K122").Value = 8994
$ws.Range("L122").Value = 8769
$ws.Range("M122").Value = -6544
$ws.Range("N122").Value = -13669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 74998
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 74998
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = ""
$ws.Range("M61").Value = 74998
$ws.Range("N61").Value = -75624
$ws.Range("H94").Value = 3312.4
$ws.Range("I94").Value = 3012.8
$ws.Range("J94").Value = 4510.8
$ws.Range("K94").Value = 3012.8
$ws.Range("L94").Value = 4510.8
$ws.Range("M94").Value = -2561.8
$ws.Range("N94").Value = -5412.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7974.5
$ws.Range("J31").Value = 16249.25
$ws.Range("L31").Value = 16249.25
$ws.Range("N31").Value = -16839.25
$ws.Range("H34").Value = 7974.5
$ws.Range("J34").Value = 16249.25
$ws.Range("L34").Value = 16249.25
$ws.Range("N34").Value = -16653.25
$ws.Range("H50").Value = 40999
$ws.Range("J50").Value = 40999
$ws.Range("L50").Value = 40999
$ws.Range("N50").Value = -42249
$ws.Range("H51").Value = 30452.273
$ws.Range("J51").Value = 30452.273
$ws.Range("L51").Value = 30452.273
$ws.Range("N51").Value = -31924.273
$ws.Range("H58").Value = 2036.625
$ws.Range("I58").Value = 2072.4
$ws.Range("K58").Value = 2072.4
$ws.Range("M58").Value = -1869.4
$ws.Range("H61").Value = 30452.273
$ws.Range("J61").Value = 30452.273
$ws.Range("L61").Value = 30452.273
$ws.Range("N61").Value = -31148.273
$ws.Range("H107").Value = 1479
$ws.Range("I107").Value = 800.1429000000001
$ws.Range("J107").Value = 1875
$ws.Range("K107").Value = 800.1429000000001
$ws.Range("L107").Value = 1875
$ws.Range("M107").Value = 1119.8571
$ws.Range("N107").Value = -5715
$ws.Range("H136").Value = 2036.625
$ws.Range("I136").Value = 2072.4
$ws.Range("K136").Value = 6217.200000000001
$ws.Range("M136").Value = -3667.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 27757.092
$ws.Range("I74").Value = 24632
$ws.Range("K74").Value = 73896
$ws.Range("M74").Value = -72835
$ws.Range("H77").Value = 27757.092
$ws.Range("I77").Value = 24632
$ws.Range("K77").Value = 221688
$ws.Range("M77").Value = -216384
$ws.Range("H112").Value = 15000
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
$ws.Range("H114").Value = 882.8333
$ws.Range("I114").Value = 366.66666
$ws.Range("J114").Value = 1399
$ws.Range("K114").Value = 1099.99998
$ws.Range("L114").Value = 4197
$ws.Range("M114").Value = 2154.00002
$ws.Range("N114").Value = -10705
$ws.Range("H117").Value = 1785.4736
$ws.Range("J117").Value = 1444.6666
$ws.Range("L117").Value = 4333.9998
$ws.Range("N117").Value = -11217.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1687
$ws.Range("I6").Value = 1549
$ws.Range("J6").Value = 1733
$ws.Range("K6").Value = 1549
$ws.Range("L6").Value = 1733
$ws.Range("M6").Value = -1436
$ws.Range("N6").Value = -1959
$ws.Range("H16").Value = 1687
$ws.Range("I16").Value = 1549
$ws.Range("J16").Value = 1733
$ws.Range("K16").Value = 1549
$ws.Range("L16").Value = 1733
$ws.Range("M16").Value = -1299
$ws.Range("N16").Value = -2233
$ws.Range("H27").Value = 5333.3335
$ws.Range("J27").Value = 5333.3335
$ws.Range("L27").Value = 5333.3335
$ws.Range("N27").Value = -5665.3335
$ws.Range("H31").Value = 1363
$ws.Range("I31").Value = 1363
$ws.Range("K31").Value = 1363
$ws.Range("M31").Value = -1071
$ws.Range("H37").Value = 1363
$ws.Range("I37").Value = 1363
$ws.Range("K37").Value = 1363
$ws.Range("M37").Value = -1086
$ws.Range("H104").Value = 33440
$ws.Range("J104").Value = 33440
$ws.Range("L104").Value = 33440
$ws.Range("N104").Value = -40428
$ws.Range("H113").Value = 2563
$ws.Range("I113").Value = 2798.6428
$ws.Range("K113").Value = 2798.6428
$ws.Range("M113").Value = -628.6428000000001
$ws.Range("H122").Value = 13810.368
$ws.Range("I122").Value = 16760.268
$ws.Range("K122").Value = 50280.804
$ws.Range("M122").Value = -47830.804

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 23700
$ws.Range("J38").Value = 23700
$ws.Range("L38").Value = 23700
$ws.Range("N38").Value = -24520
$ws.Range("H46").Value = 3122.9333
$ws.Range("I46").Value = 2235.7144
$ws.Range("J46").Value = 3899.25
$ws.Range("K46").Value = 2235.7144
$ws.Range("L46").Value = 3899.25
$ws.Range("M46").Value = -2047.7144
$ws.Range("N46").Value = -4275.25
$ws.Range("H70").Value = 45000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = ""
$ws.Range("H73").Value = 45000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = ""
$ws.Range("H136").Value = 6434
$ws.Range("I136").Value = 6266.3335
$ws.Range("J136").Value = 6601.6665
$ws.Range("K136").Value = 18799.0005
$ws.Range("L136").Value = 19804.9995
$ws.Range("M136").Value = -16249.0005
$ws.Range("N136").Value = -24904.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 36666.668
$ws.Range("J15").Value = 36666.668
$ws.Range("L15").Value = 36666.668
$ws.Range("N15").Value = -37242.668
$ws.Range("H62").Value = 16334.429
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 17422.264
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 17422.264
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -18670.264
$ws.Range("H65").Value = 16334.429
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 17422.264
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 87111.31999999999
$ws.Range("M65").Value = -26880
